$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "etc / parcel-delivery" category that used to live in column A.
# Deleting the entire column shifts B:H left into A:G, matching the new layout.
$ws.Range("A1").EntireColumn.Delete()

# Row 5 (소박스) no longer carries the duplicated "1번" entries that used to
# sit in columns D/E (now C/D after the shift) - clear them out.
$ws.Range("C5:D5").ClearContents()

# Update the two counters in row 1.
$ws.Range("A1").Value = 9
$ws.Range("B1").Value = 7

# Match the author's last selection before saving.
$ws.Range("D12").Select()
